$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the client email value (row 2, column B) to reflect new test data
$ws.Range("B2").Value = "hendi19@qh.com4"

# Update the selected/active cell in the sheet view
$ws.Range("C5").Select()
